$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 11 formatting: it should use the "last row of group" style (border),
# matching rows 4 / 6 / 8 / (new) 13. Copy formats from row 8, which already has
# that style, onto row 11 without touching its values.
$ws.Range("A8:E8").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# --- New entry: row 12 (first/only row of a new translation group), copy the
# "first row of group" formatting (no border) from row 9.
$ws.Range("A9:E9").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

$ws.Range("A12").Value = "SCRIPT/T01P01A/um2407.ssb"
$ws.Range("B12").Value = 82
$ws.Range("C12").Value = " Your team won't fail! Go for it!"
$ws.Range("D12").Value = " Ваша команда нас не подведёт!\nТолько вперёд!"
$ws.Range("E12").Value = " Âàšà ëïíàîäà îàò îå ðïäâåäæó!\nÓïìûëï âðåñæä!"
$ws.Rows.Item(12).RowHeight = 43.2

# --- Row 13: trailing blank spacer row, using the "last row of group" style
# (border), matching rows 4 / 6 / 8 / 11. Copy formats from row 11 (now fixed).
$ws.Range("A11:E11").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A13:E13").ClearContents()

$excel.CutCopyMode = $false
